$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Revert the Instance value from "Automation2" back to "Automation3"
$ws.Range("D2").Value = "Automation3"

# Revert the TestCases value from "42,43,44,45" back to "34"
$ws.Range("B2").Value = "34"
